# cap nhat nang suat - update productivity/level header values and reorder
# employee rows by base salary (Minh now first), trimming the level columns
# from 10 (up to 55) down to 7 (up to 60) and moving the "Ti le"/"Bat cap"
# columns left to follow the new last level column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (levels) ---------------------------------------------------
$ws.Range("F1").Value = 30
$ws.Range("G1").Value = 40
$ws.Range("H1").Value = 50
$ws.Range("I1").Value = 60
$ws.Range("J1").Value = "Ti le"
$ws.Range("K1").Value = "Bat cap"

# --- Row 2: Minh (moved up from row 4, flattened rates) --------------------
$ws.Range("A2").Value = "Minh"
$ws.Range("B2").Value = 4900
$ws.Range("C2").Value = 0.01
$ws.Range("D2").Value = 0.01
$ws.Range("E2").Value = 0.01
$ws.Range("F2").Value = 0.01
$ws.Range("G2").Value = 0.01
$ws.Range("H2").Value = 0.01
$ws.Range("I2").Value = 0.01

# --- Row 3: Hai (moved down from row 2) -------------------------------------
$ws.Range("A3").Value = "Hai"
$ws.Range("C3").Value = 0.011

# --- Row 4: Tien (moved down from row 3) ------------------------------------
$ws.Range("A4").Value = "Tien"
$ws.Range("B4").Value = 3000
$ws.Range("D4").Value = 0.011
$ws.Range("E4").Value = 0.012
$ws.Range("F4").Value = 0.013
$ws.Range("G4").Value = 0.014
$ws.Range("H4").Value = 0.015
$ws.Range("I4").Value = 0.016

# --- Row 5: Cuong (same position, shifted rate columns) ---------------------
$ws.Range("D5").Value = 0.011
$ws.Range("E5").Value = 0.012
$ws.Range("F5").Value = 0.013
$ws.Range("G5").Value = 0.014
$ws.Range("H5").Value = 0.015
$ws.Range("I5").Value = 0.016
$ws.Range("J5").Value = 0.7
$ws.Range("K5").Value = "*"

# --- Row 6: Duc (same position, shifted rate columns) -----------------------
$ws.Range("D6").Value = 0.011
$ws.Range("E6").Value = 0.012
$ws.Range("F6").Value = 0.013
$ws.Range("G6").Value = 0.014
$ws.Range("H6").Value = 0.015
$ws.Range("I6").Value = 0.016

# --- Clear the now-unused trailing columns (old L:N / shifted J:L) ---------
# Use Clear() (not ClearContents) so the vacated cells drop out of the sheet
# entirely (matching the shrunk A1:K12 dimension), not just lose their value.
$ws.Range("L1:N1").Clear()
$ws.Range("J2:L2").Clear()
$ws.Range("J3:L3").Clear()
$ws.Range("J4:L4").Clear()
$ws.Range("L5:N5").Clear()
$ws.Range("J6:L6").Clear()

# --- Column widths: col B narrower, cols C:I uniform width 6 ---------------
# (ColumnWidth is in characters and gets pixel-quantised on write; these are
# the character values that round-trip to the saved widths closest to the
# target 12.42578125 / 6 "best fit" pixel widths.)
$ws.Columns("B:B").ColumnWidth = 11.6
$ws.Columns("C:I").ColumnWidth = 5.14

# --- Selection cosmetics, matching the saved view ---------------------------
$ws.Range("J12").Select()
